# Auto-generated: update Price (D) and Volume(1h) (E) columns per the
# upstream cryptos-list refresh. Values are plain text in the source sheet
# (thousands separated with ".", percentages padded with spaces), so any
# replacement that Excel would otherwise parse as a genuine number is typed
# with a leading apostrophe to force literal-text entry, matching the
# original cell contents exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.956.78"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.636.09"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'212.26"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "'0.0885"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "1.639.97"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "'65.45"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "27.963.84"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "'230.95"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "'10.39"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D25").Value = "'154.94"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'1.19"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'3.41"
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").Value = "1.407.75"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("E36").Value = "  +8.69%  "
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").Value = "'0.0170"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "'0.872"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'66.88"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").Value = "1.777.60"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").Value = "'88.06"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("E51").Value = "  -0.41%  "
